# Natmi following Dr Hou advice
#
# Recomputes the LR-pair statistics for the Wnt2-Fzd3 sheet now that a new
# sending/target cluster "M2" has entered the analysis:
#   - Row 2 (Target=ECs), Row 3 (Target=FAPs) and Row 4 get refreshed
#     numeric columns E:T.
#   - Row 4's target cluster becomes the new "M2" cluster (its values are
#     recomputed for M2).
#   - A new Row 5 is appended holding the recomputed values for the
#     original "sCs" target cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($Row, $TargetCluster, $Values)

    $arr = New-Object 'object[,]' 1,17
    $arr[0,0] = $TargetCluster
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $arr[0, $i + 1] = $Values[$i]
    }
    $addr = "D" + $Row + ":T" + $Row
    $ws.Range($addr).Value = $arr
}

# Row 2: Sending=FAPs, Ligand=Wnt2, Receptor=Fzd3, Target=ECs
Set-RowData 2 "ECs" @(
    3, 1, 0.3884013333333334, 1.165204, 1, 1,
    3, 1, 0.2078313333333333, 0.623494,
    0.08621557350328635, 0.112461889302165,
    0.08072196697511112, 0.7264977027760001,
    0.08621557350328635, 0.112461889302165
)

# Row 3: Sending=FAPs, Ligand=Wnt2, Receptor=Fzd3, Target=FAPs
Set-RowData 3 "FAPs" @(
    3, 1, 0.3884013333333334, 1.165204, 1, 1,
    3, 1, 0.498127, 1.494381,
    0.206640184103479, 0.2695469573039334,
    0.1934731909693334, 1.741258718724001,
    0.206640184103479, 0.2695469573039334
)

# Row 4: Sending=FAPs, Ligand=Wnt2, Receptor=Fzd3, Target=M2 (new cluster,
# this row previously held the "sCs" data)
Set-RowData 4 "M2" @(
    3, 1, 0.3884013333333334, 1.165204, 1, 1,
    1, 0.3333333333333333, 0.016887, 0.050661,
    0.007005307459654767, 0.009137909545139137,
    0.006558933316000001, 0.05903039984400001,
    0.007005307459654767, 0.009137909545139137
)

# New Row 5: Sending=FAPs, Ligand=Wnt2, Receptor=Fzd3, Target=sCs
# (the recomputed data that used to live in row 4)
$ws.Range("A5:C5").Value = $ws.Range("A4:C4").Value()
Set-RowData 5 "sCs" @(
    3, 1, 0.3884013333333334, 1.165204, 1, 1,
    2, 1, 1.6877555, 3.375511,
    0.7001389349335798, 0.6088532438487625,
    0.6555264865406668, 3.933158919244001,
    0.7001389349335798, 0.6088532438487625
)
